# Update header labels on sheet "TwoxTwoOutTax_1-2"
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("TwoxTwoOutTax_1-2")
$ws.Activate()

# Column headers I1 / J1 swap their S0/S1 substitution-elasticity label.
# Leading apostrophe preserves the text/quote-prefix cell formatting that
# these label cells already carry (same as retyping the label in Excel).
$ws.Range("I1").Value = "'S1Tr2,1.5"
$ws.Range("J1").Value = "'S0Tr=3,1"

# Recalculated numeric results (Outtax non-0 tr_elas run)
$ws.Range("J2").Value = 1.0000000000000002
$ws.Range("K2").Value = 1.0000000000000002
$ws.Range("L2").Value = 0.87514262258757747
$ws.Range("M2").Value = 0.65301554051929678
$ws.Range("N2").Value = 0.71778700395030381
$ws.Range("O2").Value = 0.7423618925367752
$ws.Range("P2").Value = 0.54435839325680269
$ws.Range("J3").Value = 0.99999999999999989
$ws.Range("K3").Value = 0.99999999999999989
$ws.Range("L3").Value = 1.1244941080643807
$ws.Range("M3").Value = 1.3442824845817996
$ws.Range("N3").Value = 1.2804060756948938
$ws.Range("O3").Value = 1.2561259176857693
$ws.Range("P3").Value = 1.4510658620714891
$ws.Range("Q3").Value = 1.9797958971136493
$ws.Range("J4").Value = 0.99836641864624298
$ws.Range("K4").Value = 0.99893794959861781
$ws.Range("L4").Value = 0.99422149404252358
$ws.Range("M4").Value = 0.98334630994723327
$ws.Range("N4").Value = 0.98762504070267221
$ws.Range("O4").Value = 0.98786684869299268
$ws.Range("P4").Value = 0.96551294640069951
$ws.Range("J5").Value = 1.0831733772497101
$ws.Range("K5").Value = 1.0557280901465067
$ws.Range("L5").Value = 1.1552359587842187
$ws.Range("N5").Value = 1.1555980468943421
$ws.Range("O5").Value = 1.1575889833220443
$ws.Range("P5").Value = 1.2537168559460672
$ws.Range("Q5").Value = 1.5055868680147451
$ws.Range("J6").Value = 0.91682662275028926
$ws.Range("K6").Value = 0.9442719098534933
$ws.Range("L6").Value = 0.85595981732131754
$ws.Range("N6").Value = 0.85564819249958746
$ws.Range("O6").Value = 0.85393661198955839
$ws.Range("P6").Value = 0.77493692217449794
$ws.Range("Q6").Value = 0.59749254452861256
$ws.Range("J8").Value = 1.1026615242409006
$ws.Range("K8").Value = 1.1718933282765274
$ws.Range("L8").Value = 0.92664016290016515
$ws.Range("N8").Value = 0.94793633761497686
$ws.Range("O8").Value = 0.94504731852562951
$ws.Range("Q8").Value = 0.93690421845409066
$ws.Range("J9").Value = 0.72706545913678722
$ws.Range("K9").Value = 0.661136947939996
$ws.Range("L9").Value = 0.90004342873448051
$ws.Range("M9").Value = 0.88509599493431845
$ws.Range("N9").Value = 0.88903972814453536
$ws.Range("O9").Value = 0.89108439606482148
$ws.Range("P9").Value = 0.84027401686099956
$ws.Range("J10").Value = 199.67328372924806
$ws.Range("K10").Value = 199.787589545804
$ws.Range("L10").Value = 198.84429879031092
$ws.Range("M10").Value = 196.6692619893301
$ws.Range("O10").Value = 197.57336970586016
$ws.Range("P10").Value = 193.10258927643693
$ws.Range("J11").Value = 77.091007507421452
$ws.Range("K11").Value = 78.072005841314947
$ws.Range("N11").Value = 83.350496235799483
$ws.Range("O11").Value = 81.225339410207084
$ws.Range("J12").Value = 22.826613486796894
$ws.Range("K12").Value = 21.821789018740443
$ws.Range("N12").Value = 16.521274681467379
$ws.Range("O12").Value = 18.724611062026483
$ws.Range("P12").Value = 16.11595256774606
$ws.Range("J13").Value = 22.745634357202448
$ws.Range("K13").Value = 21.821789028458028
$ws.Range("N13").Value = 25.027423604599669
$ws.Range("O13").Value = 25.091438252400611
$ws.Range("P13").Value = 28.125819289701191
$ws.Range("Q13").Value = 35.028896369598066
$ws.Range("J14").Value = 77.010028377827553
$ws.Range("K14").Value = 78.072005830450337
$ws.Range("N14").Value = 74.124977375131166
$ws.Range("O14").Value = 74.038361041448425
$ws.Range("P14").Value = 69.539579788292428
$ws.Range("Q14").Value = 55.60490694634634
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 40
$ws.Range("L15").Value = 38.962246115576093
$ws.Range("N15").Value = 37.742208702045119
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 60
$ws.Range("L16").Value = 61.052960182836976
$ws.Range("M16").Value = 62.841791142654181
$ws.Range("N16").Value = 62.331212543599889
$ws.Range("O16").Value = 62.135132465669287
$ws.Range("P16").Value = 63.678287301078129
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 60
$ws.Range("L17").Value = 58.606334458651354
$ws.Range("M17").Value = 56.286003007094578
$ws.Range("N17").Value = 56.942195509126392
$ws.Range("O17").Value = 57.195490828323443
$ws.Range("P17").Value = 55.221484067444941
$ws.Range("Q17").Value = 50.510257216772395
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 40
$ws.Range("L18").Value = 41.41413631445721
$ws.Range("M18").Value = 43.86230905047411
$ws.Range("N18").Value = 43.157765918714325
$ws.Range("O18").Value = 42.888411696208422
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 100
$ws.Range("L20").Value = 93.038917910466012
$ws.Range("N20").Value = 93.02434064229395
$ws.Range("O20").Value = 92.944309968961264
$ws.Range("P20").Value = 89.310036642749893
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 100
$ws.Range("L21").Value = 108.08696225125264
$ws.Range("N21").Value = 108.10664295323191
$ws.Range("O21").Value = 108.21493002930814
$ws.Range("P21").Value = 113.59698981509561
$ws.Range("Q21").Value = 129.37005259843758
$ws.Range("J22").Value = 199.67328372924806
$ws.Range("K22").Value = 199.787589545804
$ws.Range("L22").Value = 198.84429879031092
$ws.Range("M22").Value = 196.6692619893301
$ws.Range("O22").Value = 197.57336970586016
$ws.Range("P22").Value = 193.10258927643693
$ws.Range("J23").Value = 0.99836641864624032
$ws.Range("K23").Value = 0.99893794772902
$ws.Range("L23").Value = 0.99422149395155457
$ws.Range("M23").Value = 0.98334630994665051
$ws.Range("O23").Value = 0.98786684852930084
$ws.Range("P23").Value = 0.96551294638218466

Write-Output "edit complete"
